# Applies the diff described in the commit:
#  - swap match details (columns F:V) between rows 4 & 5
#  - swap match details (columns F:V) between rows 24 & 25
#  - swap match details (columns F:V) between rows 86 & 87
#  - append three new match rows (99, 100, 101)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchDetails {
    param($rowA, $rowB)

    # Stash row A's F:V values (literal Value2, so formatting/types survive)
    $vals = @{}
    for ($col = 6; $col -le 22; $col++) {
        $vals[$col] = $ws.Cells.Item($rowA, $col).Value2
    }

    # Row A <- Row B
    for ($col = 6; $col -le 22; $col++) {
        $ws.Cells.Item($rowA, $col).Value = $ws.Cells.Item($rowB, $col).Value2
    }

    # Row B <- stashed Row A
    for ($col = 6; $col -le 22; $col++) {
        $ws.Cells.Item($rowB, $col).Value = $vals[$col]
    }
}

Swap-MatchDetails 4 5
Swap-MatchDetails 24 25
Swap-MatchDetails 86 87

# Append the three new rows, copying formatting (styles) from the last
# existing data row (98) so the index/date columns keep their look.
$ws.Range("A98:V98").Copy()
$ws.Range("A99:V101").PasteSpecial(-4122)

function Set-MatchRow {
    param(
        $row, $idx, $matchDate,
        $home, $homeGoals, $away, $awayGoals,
        $homeOpen, $homeOpenDt, $homeClose, $homeCloseDt,
        $drawOpen, $drawOpenDt, $drawClose, $drawCloseDt,
        $awayOpen, $awayOpenDt, $awayClose, $awayCloseDt,
        $url
    )

    $ws.Cells.Item($row, 1).Value = $idx
    $ws.Cells.Item($row, 2).Value = "israel"
    $ws.Cells.Item($row, 3).Value = "ligat-ha-al"
    $ws.Cells.Item($row, 4).Value = "2023-2024"
    $ws.Cells.Item($row, 5).Value = $matchDate
    $ws.Cells.Item($row, 6).Value = $home
    $ws.Cells.Item($row, 7).Value = $homeGoals
    $ws.Cells.Item($row, 8).Value = $away
    $ws.Cells.Item($row, 9).Value = $awayGoals
    $ws.Cells.Item($row, 10).Value = $homeOpen
    $ws.Cells.Item($row, 11).Value = $homeOpenDt
    $ws.Cells.Item($row, 12).Value = $homeClose
    $ws.Cells.Item($row, 13).Value = $homeCloseDt
    $ws.Cells.Item($row, 14).Value = $drawOpen
    $ws.Cells.Item($row, 15).Value = $drawOpenDt
    $ws.Cells.Item($row, 16).Value = $drawClose
    $ws.Cells.Item($row, 17).Value = $drawCloseDt
    $ws.Cells.Item($row, 18).Value = $awayOpen
    $ws.Cells.Item($row, 19).Value = $awayOpenDt
    $ws.Cells.Item($row, 20).Value = $awayClose
    $ws.Cells.Item($row, 21).Value = $awayCloseDt
    $ws.Cells.Item($row, 22).Value = $url
}

Set-MatchRow 99 98 45297.58333333334 `
    "H. Beer Sheva" 4 "Maccabi Bnei Raina" 2 `
    1.58 "02/01/2024 19:12" 1.48 "06/01/2024 13:58" `
    3.87 "02/01/2024 19:12" 4.26 "06/01/2024 13:58" `
    5.34 "02/01/2024 19:12" 7.17 "06/01/2024 13:58" `
    "https://www.betexplorer.com/football/israel/ligat-ha-al/h-beer-sheva-maccabi-bnei-raina/QNBw8rzd/"

Set-MatchRow 100 99 45297.70833333334 `
    "Sakhnin" 0 "Hapoel Petah Tikva" 0 `
    2.05 "03/01/2024 19:42" 1.99 "06/01/2024 16:52" `
    3.27 "03/01/2024 19:42" 3.41 "06/01/2024 16:52" `
    3.52 "03/01/2024 19:42" 4 "06/01/2024 16:52" `
    "https://www.betexplorer.com/football/israel/ligat-ha-al/sakhnin-hapoel-petah-tikva/6JlQVoS2/"

Set-MatchRow 101 100 45297.77083333334 `
    "Netanya" 2 "Hapoel Tel Aviv" 1 `
    1.94 "05/01/2024 14:29" 1.86 "06/01/2024 18:28" `
    3.4 "05/01/2024 14:29" 3.81 "06/01/2024 18:28" `
    4.22 "05/01/2024 14:29" 4.03 "06/01/2024 18:28" `
    "https://www.betexplorer.com/football/israel/ligat-ha-al/netanya-hapoel-tel-aviv/4EAs72k2/"

# Keep the declared dimension in sync with the new used range.
$ws.Range("A1:V101").Select()
